# CIERRE 20 NOV 2021
# Advance the payroll workbook from "SEMANA 46 (08-14 NOV 2021)" to
# "SEMANA 47 (15-21 NOV 2021)" and update the week's pay figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 holds the week-range banner; H9/B27/H27/B43/H43/B60 all reference it
# (directly or transitively) via formulas, so updating this one cell
# cascades the new week label everywhere it is shown.
$ws.Range("B9").Value = "SEMANA   47  DEL    15      Al    21   DE   NOVIEMBRE          2021"

# Weekly pay entries that changed for this close.
$ws.Range("K4").Value = 0
$ws.Range("K21").Value = 1806
$ws.Range("K39").Value = 1250

# Restore the on-screen scroll position / selection as last saved.
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("H60").Select()
